$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: replace the long concatenated "everything" text with just the abstract
$ws.Range("D2").Value = "Our current understanding of the spectrum of TB and COVID-19 lesions in the human lung is limited by a reliance on low-resolution imaging platforms that cannot provide accurate 3D representations of lesion types within the context of the whole lung. To characterize TB and COVID-19 lesions in 3D, we applied micro/nanocomputed tomography to surgically resected, postmortem, and paraffin-embedded human lung tissue. We define a spectrum of TB pathologies, including cavitary lesions, calcium deposits outside and inside necrotic granulomas and mycetomas, and vascular rearrangement. We identified an unusual spatial arrangement of vasculature within an entire COVID-19 lobe, and 3D segmentation of blood vessels revealed microangiopathy associated with hemorrhage. Notably, segmentation of pathological anomalies reveals hidden pathological structures that might otherwise be disregarded, demonstrating a powerful method to visualize pathologies in 3D in TB lung tissue and whole COVID-19 lobes. These findings provide unexpected new insight into the spatial organization of the spectrum of TB and COVID-19 lesions within the framework of the entire lung."

# X2: reorder the National Institutes of Health award_id entries (same line-wrap structure)
$x2Line1 = "c(funder = `"https://openalex.org/F4320306137`", funder_display_name = `"Bill and Melinda Gates Foundation`", award_id = `"OPP1116944`", funder = `"https://openalex.org/F4320306137`", funder_display_name = `"Bill and Melinda Gates Foundation`", award_id = `"OPP1130017`", funder = `"https://openalex.org/F4320306168`", funder_display_name = `"Cystic Fibrosis Foundation`", award_id = `"ROWE19R0`", funder = `"https://openalex.org/F4320307853`", funder_display_name = `"CRDF Global`", funder = `"https://openalex.org/F4320320671`", "
$x2Line2 = "funder_display_name = `"National Research Foundation`", funder = `"https://openalex.org/F4320320672`", funder_display_name = `"South African Medical Research Council`", funder = `"https://openalex.org/F4320332161`", funder_display_name = `"National Institutes of Health`", award_id = `"U01 HL152978`", funder = `"https://openalex.org/F4320332161`", funder_display_name = `"National Institutes of Health`", award_id = `"5F31HL146083‐02`", funder = `"https://openalex.org/F4320332161`", funder_display_name = `"National Institutes of Health`", "
$x2Line3 = "award_id = `"R21A127182`", funder = `"https://openalex.org/F4320332161`", funder_display_name = `"National Institutes of Health`", award_id = `"P30 DK072482`", funder = `"https://openalex.org/F4320332161`", funder_display_name = `"National Institutes of Health`", award_id = `"3T32GM008361‐30S1`", funder = `"https://openalex.org/F4320332161`", funder_display_name = `"National Institutes of Health`", award_id = `"R33AI138280`", funder = `"https://openalex.org/F4320332161`", funder_display_name = `"National Institutes of Health`", "
$x2Line4 = "award_id = `"2T32HL105346‐11A1`", funder = `"https://openalex.org/F4320332161`", funder_display_name = `"National Institutes of Health`", award_id = `"R01Al111940`", funder = `"https://openalex.org/F4320332161`", funder_display_name = `"National Institutes of Health`", award_id = `"R35 HL135816‐04S1`", funder = `"https://openalex.org/F4320334435`", funder_display_name = `"Center for AIDS Research, University of Alabama at Birmingham`")"

$ws.Range("X2").Value = $x2Line1 + "`n" + $x2Line2 + "`n" + $x2Line3 + "`n" + $x2Line4

# AF2: reorder related_works list
$ws.Range("AF2").Value = "c(`"https://openalex.org/W4206669628`", `"https://openalex.org/W3036314732`", `"https://openalex.org/W4205317059`", `"https://openalex.org/W3009669391`", `"https://openalex.org/W4382894326`", `"https://openalex.org/W3198183218`", `"https://openalex.org/W3081785542`", `"https://openalex.org/W3176864053`", `"https://openalex.org/W3084498529`", `"https://openalex.org/W3171943759`")"
